# Updates cryptos list figures (Price + Volume(1h)) per the Jun 17 2023 refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.456.55"
$ws.Range("E2").Value = "  +3.48%  "
$ws.Range("D3").Value = "1.731.74"
$ws.Range("E3").Value = "  +3.93%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.33%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.66"
$ws.Range("E5").Value = "  +2.95%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  +0.33%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4785"
$ws.Range("E7").Value = "  +3.26%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2661"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06219"
$ws.Range("E9").Value = "  +1.27%  "
$ws.Range("D10").Value = "1.736.30"
$ws.Range("E10").Value = "  +4.20%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07122"
$ws.Range("E11").Value = "  +2.68%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.65"
$ws.Range("E12").Value = "  +5.23%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6125"
$ws.Range("E13").Value = "  +7.01%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.517"
$ws.Range("E14").Value = "  +4.02%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "76.81"
$ws.Range("E15").Value = "  +2.20%  "
$ws.Range("E16").Value = "  +0.18%  "
$ws.Range("D17").Value = "26.490.48"
$ws.Range("E17").Value = "  +3.66%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.001"
$ws.Range("E18").Value = "  +0.25%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000006889"
$ws.Range("E19").Value = "  +2.55%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.69"
$ws.Range("E20").Value = "  +2.57%  "
$ws.Range("D21").Value = "1.958.76"
$ws.Range("E21").Value = "  +4.54%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.561"
$ws.Range("E22").Value = "  +3.04%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.872"
$ws.Range("E23").Value = "  +2.89%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.326"
$ws.Range("E24").Value = "  +1.85%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "135.94"
$ws.Range("E25").Value = "  +1.34%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.33"
$ws.Range("E26").Value = "  +2.57%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.791"
$ws.Range("E27").Value = "  +4.19%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.398"
$ws.Range("E28").Value = "  +1.92%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "106.58"
$ws.Range("E29").Value = "  +2.30%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.975"
$ws.Range("E30").Value = "  +0.62%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.07973"
$ws.Range("E31").Value = "  +3.71%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.705"
$ws.Range("E32").Value = "  +2.97%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04541"
$ws.Range("E33").Value = "  +4.32%  "
$ws.Range("E34").Value = "  +0.43%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6348"
$ws.Range("E35").Value = "  +4.59%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9890"
$ws.Range("E36").Value = "  +5.16%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9310"
$ws.Range("E37").Value = "  +1.87%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "109.92"
$ws.Range("E38").Value = "  +1.75%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.975"
$ws.Range("E39").Value = "  +7.52%  "
$ws.Range("E40").Value = "  +0.61%  "
$ws.Range("E41").Value = "  +0.70%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.01503"
$ws.Range("E42").Value = "  +3.50%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.695"
$ws.Range("E43").Value = "  +13.34%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.3890"
$ws.Range("E44").Value = "  +4.57%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "6.896"
$ws.Range("E45").Value = "  +12.68%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1189"
$ws.Range("E46").Value = "  +7.10%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05333"
$ws.Range("E47").Value = "  +1.37%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.881"
$ws.Range("E48").Value = "  +4.24%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "30.74"
$ws.Range("E49").Value = "  +0.53%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.265"
$ws.Range("E50").Value = "  +5.29%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3414"
$ws.Range("E51").Value = "  +2.87%  "
